$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rewrite the data rows (2-10) in their new (alphabetically-by-Spanish-verb)
# order. One verb/word pair is also new: "Leli" now pairs with the new word
# "ojear" (previously it was paired with "mirar", which is now used by
# "Kintu" instead).

$ws.Range("A2").Value = "Püra"
$ws.Range("B2").Value = "ascender"

$ws.Range("A3").Value = "Pi"
$ws.Range("B3").Value = "decir"

$ws.Range("A4").Value = "Kon"
$ws.Range("B4").Value = "entrar"

$ws.Range("A5").Value = "Kintu"
$ws.Range("B5").Value = "mirar"

$ws.Range("A6").Value = "Leli"
$ws.Range("B6").Value = "ojear"

$ws.Range("A7").Value = "Pepi"
$ws.Range("B7").Value = "poder"

$ws.Range("A8").Value = "Wüño"
$ws.Range("B8").Value = "regresar"

$ws.Range("A9").Value = "Tripa"
$ws.Range("B9").Value = "salir"

$ws.Range("A10").Value = "Rüngkü"
$ws.Range("B10").Value = "saltar"

# Column A keeps the original (no explicit style) look for the words that
# never carried the extra style, regardless of which row they ended up in.
$ws.Range("A2").ClearFormats()
$ws.Range("A5").ClearFormats()
$ws.Range("A7").ClearFormats()
$ws.Range("A8").ClearFormats()
$ws.Range("A10").ClearFormats()

# Move the active selection, as recorded in the saved workbook.
$ws.Range("E11").Select()
